$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

# Locate the paragraph that currently contains "Streamlit nous offre un frontend interactif qui"
$targetIndex = -1
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "Streamlit nous offre un frontend interactif qui*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate target paragraph"
}

$pMain = $d.Paragraphs.Item($targetIndex)
$pNext = $d.Paragraphs.Item($targetIndex + 1)

# The paragraph right after $pNext should be empty; if so, extend the
# replacement range to swallow it (it disappears in the target layout).
$endOfRange = $pMain.Range.End
if ($pNext.Range.Text.Trim().Length -eq 0) {
    $endOfRange = $pNext.Range.End
}

$replaceRange = $d.Range($pMain.Range.Start, $endOfRange)

$xmlMain = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:t>L''outil est accessible aux adresses suivantes</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">URL Railway: </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:instrText xml:space="preserve"> HYPERLINK "https://backend-scoring.up.railway.app" </w:instrText></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rStyle w:val="3"/><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:t>h</w:t></w:r><w:r><w:rPr><w:rStyle w:val="3"/><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:t>ttps://backend-scoring.up.railway.app</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:t xml:space="preserve">URL Streamlit: </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:instrText xml:space="preserve"> HYPERLINK "https://projet-7-2fdf9ahp4nvj6yinucaj9q.streamlit.app/" </w:instrText></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rStyle w:val="3"/><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:t>https://projet-7-2fdf9ahp4nvj6yinucaj9q.streamlit.app/</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/styles.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.styles+xml"><pkg:xmlData><?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:styles xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:sl="http://schemas.openxmlformats.org/schemaLibrary/2006/main" xmlns:wpsCustomData="http://www.wps.cn/officeDocument/2013/wpsCustomData" mc:Ignorable="w14"><w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="SimSun" w:cs="Times New Roman"/></w:rPr></w:rPrDefault><w:pPrDefault/></w:docDefaults><w:latentStyles w:count="260" w:defQFormat="0" w:defUnhideWhenUsed="1" w:defSemiHidden="1" w:defUIPriority="99" w:defLockedState="0"><w:lsdException w:qFormat="1" w:unhideWhenUsed="0" w:uiPriority="0" w:semiHidden="0" w:name="Normal"/><w:lsdException w:qFormat="1" w:unhideWhenUsed="0" w:uiPriority="9" w:semiHidden="0" w:name="heading 1"/><w:lsdException w:qFormat="1" w:uiPriority="9" w:name="heading 2"/><w:lsdException w:qFormat="1" w:uiPriority="9" w:name="heading 3"/><w:lsdException w:qFormat="1" w:uiPriority="9" w:name="heading 4"/><w:lsdException w:qFormat="1" w:uiPriority="9" w:name="heading 5"/><w:lsdException w:qFormat="1" w:uiPriority="9" w:name="heading 6"/><w:lsdException w:qFormat="1" w:uiPriority="9" w:name="heading 7"/><w:lsdException w:qFormat="1" w:uiPriority="9" w:name="heading 8"/><w:lsdException w:qFormat="1" w:uiPriority="9" w:name="heading 9"/><w:lsdException w:uiPriority="99" w:name="index 1"/><w:lsdException w:uiPriority="99" w:name="index 2"/><w:lsdException w:uiPriority="99" w:name="index 3"/><w:lsdException w:uiPriority="99" w:name="index 4"/><w:lsdException w:uiPriority="99" w:name="index 5"/><w:lsdException w:uiPriority="99" w:name="index 6"/><w:lsdException w:uiPriority="99" w:name="index 7"/><w:lsdException w:uiPriority="99" w:name="index 8"/><w:lsdException w:uiPriority="99" w:name="index 9"/><w:lsdException w:uiPriority="39" w:name="toc 1"/><w:lsdException w:uiPriority="39" w:name="toc 2"/><w:lsdException w:uiPriority="39" w:name="toc 3"/><w:lsdException w:uiPriority="39" w:name="toc 4"/><w:lsdException w:uiPriority="39" w:name="toc 5"/><w:lsdException w:uiPriority="39" w:name="toc 6"/><w:lsdException w:uiPriority="39" w:name="toc 7"/><w:lsdException w:uiPriority="39" w:name="toc 8"/><w:lsdException w:uiPriority="39" w:name="toc 9"/><w:lsdException w:uiPriority="99" w:name="Normal Indent"/><w:lsdException w:uiPriority="99" w:name="footnote text"/><w:lsdException w:uiPriority="99" w:name="annotation text"/><w:lsdException w:qFormat="1" w:uiPriority="99" w:semiHidden="0" w:name="header"/><w:lsdException w:qFormat="1" w:uiPriority="99" w:semiHidden="0" w:name="footer"/><w:lsdException w:uiPriority="99" w:name="index heading"/><w:lsdException w:qFormat="1" w:uiPriority="35" w:name="caption"/><w:lsdException w:uiPriority="99" w:name="table of figures"/><w:lsdException w:uiPriority="99" w:name="envelope address"/><w:lsdException w:uiPriority="99" w:name="envelope return"/><w:lsdException w:uiPriority="99" w:name="footnote reference"/><w:lsdException w:uiPriority="99" w:name="annotation reference"/><w:lsdException w:uiPriority="99" w:name="line number"/><w:lsdException w:uiPriority="99" w:name="page number"/><w:lsdException w:uiPriority="99" w:name="endnote reference"/><w:lsdException w:uiPriority="99" w:name="endnote text"/><w:lsdException w:uiPriority="99" w:name="table of authorities"/><w:lsdException w:uiPriority="99" w:name="macro"/><w:lsdException w:uiPriority="99" w:name="toa heading"/><w:lsdException w:uiPriority="99" w:name="List"/><w:lsdException w:uiPriority="99" w:name="List Bullet"/><w:lsdException w:uiPriority="99" w:name="List Number"/><w:lsdException w:uiPriority="99" w:name="List 2"/><w:lsdException w:uiPriority="99" w:name="List 3"/><w:lsdException w:uiPriority="99" w:name="List 4"/><w:lsdException w:uiPriority="99" w:name="List 5"/><w:lsdException w:uiPriority="99" w:name="List Bullet 2"/><w:lsdException w:uiPriority="99" w:name="List Bullet 3"/><w:lsdException w:uiPriority="99" w:name="List Bullet 4"/><w:lsdException w:uiPriority="99" w:name="List Bullet 5"/><w:lsdException w:uiPriority="99" w:name="List Number 2"/><w:lsdException w:uiPriority="99" w:name="List Number 3"/><w:lsdException w:uiPriority="99" w:name="List Number 4"/><w:lsdException w:uiPriority="99" w:name="List Number 5"/><w:lsdException w:qFormat="1" w:unhideWhenUsed="0" w:uiPriority="10" w:semiHidden="0" w:name="Title"/><w:lsdException w:uiPriority="99" w:name="Closing"/><w:lsdException w:uiPriority="99" w:name="Signature"/><w:lsdException w:uiPriority="1" w:name="Default Paragraph Font"/><w:lsdException w:uiPriority="99" w:name="Body Text"/><w:lsdException w:uiPriority="99" w:name="Body Text Indent"/><w:lsdException w:uiPriority="99" w:name="List Continue"/><w:lsdException w:uiPriority="99" w:name="List Continue 2"/><w:lsdException w:uiPriority="99" w:name="List Continue 3"/><w:lsdException w:uiPriority="99" w:name="List Continue 4"/><w:lsdException w:uiPriority="99" w:name="List Continue 5"/><w:lsdException w:uiPriority="99" w:name="Message Header"/><w:lsdException w:qFormat="1" w:unhideWhenUsed="0" w:uiPriority="11" w:semiHidden="0" w:name="Subtitle"/><w:lsdException w:uiPriority="99" w:name="Salutation"/><w:lsdException w:uiPriority="99" w:name="Date"/><w:lsdException w:uiPriority="99" w:name="Body Text First Indent"/><w:lsdException w:uiPriority="99" w:name="Body Text First Indent 2"/><w:lsdException w:uiPriority="99" w:name="Note Heading"/><w:lsdException w:uiPriority="99" w:name="Body Text 2"/><w:lsdException w:uiPriority="99" w:name="Body Text 3"/><w:lsdException w:uiPriority="99" w:name="Body Text Indent 2"/><w:lsdException w:uiPriority="99" w:name="Body Text Indent 3"/><w:lsdException w:uiPriority="99" w:name="Block Text"/><w:lsdException w:qFormat="1" w:uiPriority="99" w:semiHidden="0" w:name="Hyperlink"/><w:lsdException w:uiPriority="99" w:name="FollowedHyperlink"/><w:lsdException w:qFormat="1" w:unhideWhenUsed="0" w:uiPriority="22" w:semiHidden="0" w:name="Strong"/><w:lsdException w:qFormat="1" w:unhideWhenUsed="0" w:uiPriority="20" w:semiHidden="0" w:name="Emphasis"/><w:lsdException w:uiPriority="99" w:name="Document Map"/><w:lsdException w:uiPriority="99" w:name="Plain Text"/><w:lsdException w:uiPriority="99" w:name="E-mail Signature"/><w:lsdException w:uiPriority="99" w:name="Normal (Web)"/><w:lsdException w:uiPriority="99" w:name="HTML Acronym"/><w:lsdException w:uiPriority="99" w:name="HTML Address"/><w:lsdException w:uiPriority="99" w:name="HTML Cite"/><w:lsdException w:uiPriority="99" w:name="HTML Code"/><w:lsdException w:uiPriority="99" w:name="HTML Definition"/><w:lsdException w:uiPriority="99" w:name="HTML Keyboard"/><w:lsdException w:uiPriority="99" w:name="HTML Preformatted"/><w:lsdException w:uiPriority="99" w:name="HTML Sample"/><w:lsdException w:uiPriority="99" w:name="HTML Typewriter"/><w:lsdException w:uiPriority="99" w:name="HTML Variable"/><w:lsdException w:uiPriority="99" w:name="Normal Table"/><w:lsdException w:uiPriority="99" w:name="annotation subject"/><w:lsdException w:uiPriority="99" w:name="Balloon Text"/></w:latentStyles><w:style w:type="paragraph" w:default="1" w:styleId="1"><w:name w:val="Normal"/><w:qFormat/><w:uiPriority w:val="0"/><w:pPr><w:spacing w:after="160" w:line="259" w:lineRule="auto"/></w:pPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:eastAsiaTheme="minorHAnsi" w:cstheme="minorBidi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="fr-FR" w:eastAsia="en-US" w:bidi="ar-SA"/></w:rPr></w:style><w:style w:type="character" w:default="1" w:styleId="2"><w:name w:val="Default Paragraph Font"/><w:semiHidden/><w:unhideWhenUsed/><w:uiPriority w:val="1"/></w:style><w:style w:type="table" w:default="1" w:styleId="6"><w:name w:val="Normal Table"/><w:semiHidden/><w:unhideWhenUsed/><w:uiPriority w:val="99"/><w:tblPr><w:tblCellMar><w:top w:w="0" w:type="dxa"/><w:left w:w="108" w:type="dxa"/><w:bottom w:w="0" w:type="dxa"/><w:right w:w="108" w:type="dxa"/></w:tblCellMar></w:tblPr></w:style><w:style w:type="character" w:styleId="3"><w:name w:val="Hyperlink"/><w:basedOn w:val="2"/><w:unhideWhenUsed/><w:qFormat/><w:uiPriority w:val="99"/><w:rPr><w:color w:val="0563C1" w:themeColor="hyperlink"/><w:u w:val="single"/><w14:textFill><w14:solidFill><w14:schemeClr w14:val="hlink"/></w14:solidFill></w14:textFill></w:rPr></w:style><w:style w:type="paragraph" w:styleId="4"><w:name w:val="footer"/><w:basedOn w:val="1"/><w:link w:val="9"/><w:unhideWhenUsed/><w:qFormat/><w:uiPriority w:val="99"/><w:pPr><w:tabs><w:tab w:val="center" w:pos="4536"/><w:tab w:val="right" w:pos="9072"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:style><w:style w:type="paragraph" w:styleId="5"><w:name w:val="header"/><w:basedOn w:val="1"/><w:link w:val="8"/><w:unhideWhenUsed/><w:qFormat/><w:uiPriority w:val="99"/><w:pPr><w:tabs><w:tab w:val="center" w:pos="4536"/><w:tab w:val="right" w:pos="9072"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr></w:style><w:style w:type="character" w:customStyle="1" w:styleId="7"><w:name w:val="Unresolved Mention"/><w:basedOn w:val="2"/><w:semiHidden/><w:unhideWhenUsed/><w:qFormat/><w:uiPriority w:val="99"/><w:rPr><w:color w:val="605E5C"/><w:shd w:val="clear" w:color="auto" w:fill="E1DFDD"/></w:rPr></w:style><w:style w:type="character" w:customStyle="1" w:styleId="8"><w:name w:val="En-tête Car"/><w:basedOn w:val="2"/><w:link w:val="5"/><w:qFormat/><w:uiPriority w:val="99"/></w:style><w:style w:type="character" w:customStyle="1" w:styleId="9"><w:name w:val="Pied de page Car"/><w:basedOn w:val="2"/><w:link w:val="4"/><w:qFormat/><w:uiPriority w:val="99"/></w:style></w:styles></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/styles" Target="styles.xml"/></Relationships></pkg:xmlData></pkg:part></pkg:package>'
$replaceRange.InsertXML($xmlMain)

# Now locate the paragraph holding the screenshot image (first paragraph
# from the end that contains an inline shape) and add the _GoBack bookmark
# back at its very start, matching the original document structure.
$imgParaIndex = -1
$paraCount2 = $d.Paragraphs.Count
for ($i = $paraCount2; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $imgParaIndex = $i
        break
    }
}

if ($imgParaIndex -eq -1) {
    throw "Could not locate image paragraph"
}

$pImg = $d.Paragraphs.Item($imgParaIndex)
$insPoint = $d.Range($pImg.Range.Start, $pImg.Range.Start)
$xmlBookmark = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint.InsertXML($xmlBookmark)

Write-Host "Edit applied successfully"
